# "Generate Report for Handoff"
#
# Updates the localization-status report:
#  - Status moves from "Handed back: in sync with en-US" to "In Translation"
#    (reflected on the Overview sheet as well as each language sheet).
#  - Refreshes the "Latest Handoff Datetime" timestamps.
#  - Records a new "Error Detail" warning on each language sheet noting
#    that the handback file version is stale.
#  - Narrows a couple of over-wide columns that were sized for the old,
#    longer status text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "In Translation"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a072feab344118faf8bc8e6a6507da2c56498f26/e2e/c646d137-6169-4650-991f-2d337c5289f5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a682d9e94f201a4b72085823d21499a1f85b50b0/e2e/c646d137-6169-4650-991f-2d337c5289f5.md."

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2017-02-09 13:59:50"

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2017-02-09 13:59:32"
$zhcn.Range("R2").Value = $errorDetail

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$zhcn.Columns.Item(18).ColumnWidth = 39.1666666666

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2017-02-09 13:59:50"
$dede.Range("R2").Value = $errorDetail

$dede.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(18).ColumnWidth = 39.1666666666
